# Scheduled market-data refresh: update computed price/profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ -- columns H:N) on each class/job sheet with freshly
# fetched values. Only specific rows/cells are refreshed; everything
# else in the workbook is left untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 226.83333
$ws.Range("I9").Value = 172.2
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 172.2
$ws.Range("L9").Value = 500
$ws.Range("M9").Value = -3.199999999999989
$ws.Range("N9").Value = -838
$ws.Range("H33").Value = 2334.0908
$ws.Range("I33").Value = 1630.5555
$ws.Range("K33").Value = 1630.5555
$ws.Range("M33").Value = -1401.5555
$ws.Range("H40").Value = 3666.6667
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3666.6667
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = 3666.6667
$ws.Range("N40").Value = -4016.6667
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("N44").Value = 0
$ws.Range("H70").Value = 1699
$ws.Range("I70").Value = 1168
$ws.Range("J70").Value = 2097.25
$ws.Range("K70").Value = 3504
$ws.Range("L70").Value = 6291.75
$ws.Range("M70").Value = -3234
$ws.Range("N70").Value = -6831.75
$ws.Range("H73").Value = 1699
$ws.Range("I73").Value = 1168
$ws.Range("J73").Value = 2097.25
$ws.Range("K73").Value = 3504
$ws.Range("L73").Value = 6291.75
$ws.Range("M73").Value = -2568
$ws.Range("N73").Value = -8163.75
$ws.Range("H98").Value = 2333.3333
$ws.Range("I98").Value = 1600
$ws.Range("J98").Value = 2480
$ws.Range("K98").Value = 1600
$ws.Range("L98").Value = 2480
$ws.Range("M98").Value = -102
$ws.Range("N98").Value = -5476
$ws.Range("H122").Value = 2333.3333
$ws.Range("I122").Value = 1600
$ws.Range("J122").Value = 2480
$ws.Range("K122").Value = 4800
$ws.Range("L122").Value = 7440
$ws.Range("M122").Value = -2350
$ws.Range("N122").Value = -12340
$ws.Range("H135").Value = 1044.1904
$ws.Range("I135").Value = 1211.5294
$ws.Range("J135").Value = 333
$ws.Range("K135").Value = 10903.7646
$ws.Range("L135").Value = 2997
$ws.Range("M135").Value = -8368.764599999999
$ws.Range("N135").Value = -8067
$ws.Range("H137").Value = 1697.3125
$ws.Range("I137").Value = 940.4
$ws.Range("J137").Value = 2041.3636
$ws.Range("K137").Value = 2821.2
$ws.Range("L137").Value = 6124.0908
$ws.Range("M137").Value = -271.1999999999998
$ws.Range("N137").Value = -11224.0908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10419035
$ws.Range("I61").Value = 15153445
$ws.Range("K61").Value = 15153445
$ws.Range("M61").Value = -15153233
$ws.Range("H74").Value = 882.2857
$ws.Range("I74").Value = 917.4583
$ws.Range("J74").Value = 835.3889
$ws.Range("K74").Value = 917.4583
$ws.Range("L74").Value = 835.3889
$ws.Range("M74").Value = -43.45830000000001
$ws.Range("N74").Value = -2583.3889
$ws.Range("H77").Value = 882.2857
$ws.Range("I77").Value = 917.4583
$ws.Range("J77").Value = 835.3889
$ws.Range("K77").Value = 4587.2915
$ws.Range("L77").Value = 4176.944500000001
$ws.Range("M77").Value = -219.2915000000003
$ws.Range("N77").Value = -12912.9445
$ws.Range("H122").Value = 2709.8572
$ws.Range("I122").Value = 3656.3333
$ws.Range("K122").Value = 10968.9999
$ws.Range("M122").Value = -8518.999899999999
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("N123").Value = 0
$ws.Range("H132").Value = 9236.235000000001
$ws.Range("I132").Value = 8847.385
$ws.Range("J132").Value = 10500
$ws.Range("K132").Value = 26542.155
$ws.Range("L132").Value = 31500
$ws.Range("M132").Value = -24012.155
$ws.Range("N132").Value = -36560
$ws.Range("H136").Value = 10419035
$ws.Range("I136").Value = 15153445
$ws.Range("K136").Value = 45460335
$ws.Range("M136").Value = -45457785

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 5622
$ws.Range("J6").Value = 5622
$ws.Range("L6").Value = 5622
$ws.Range("N6").Value = -5848
$ws.Range("H13").Value = 42000
$ws.Range("J13").Value = 42000
$ws.Range("L13").Value = 42000
$ws.Range("N13").Value = -42336
$ws.Range("H22").Value = 786.7222
$ws.Range("I22").Value = 772.5625
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 772.5625
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -599.5625
$ws.Range("N22").Value = -1246
$ws.Range("H50").Value = 38340
$ws.Range("J50").Value = 38340
$ws.Range("L50").Value = 38340
$ws.Range("N50").Value = -39488
$ws.Range("H94").Value = 818.1
$ws.Range("I94").Value = 845.2
$ws.Range("J94").Value = 791
$ws.Range("K94").Value = 845.2
$ws.Range("L94").Value = 791
$ws.Range("M94").Value = -394.2
$ws.Range("N94").Value = -1693
$ws.Range("H107").Value = 3000
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 0
$ws.Range("M107").Value = 3000
$ws.Range("N107").Value = -6840
$ws.Range("H115").Value = 77349.42999999999
$ws.Range("J115").Value = 77349.42999999999
$ws.Range("L115").Value = 77349.42999999999
$ws.Range("N115").Value = -80483.42999999999
$ws.Range("H119").Value = 21583.25
$ws.Range("J119").Value = 21583.25
$ws.Range("L119").Value = 21583.25
$ws.Range("N119").Value = -31259.25
$ws.Range("H134").Value = 3670.1538
$ws.Range("I134").Value = 3564.7273
$ws.Range("J134").Value = 4250
$ws.Range("K134").Value = 10694.1819
$ws.Range("L134").Value = 12750
$ws.Range("M134").Value = -8159.1819
$ws.Range("N134").Value = -17820
$ws.Range("H138").Value = 60780
$ws.Range("J138").Value = 60780
$ws.Range("L138").Value = 60780
$ws.Range("N138").Value = -71060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H114").Value = 54684
$ws.Range("J114").Value = 54684
$ws.Range("L114").Value = 54684
$ws.Range("N114").Value = -63362
$ws.Range("H117").Value = 50428
$ws.Range("J117").Value = 50428
$ws.Range("L117").Value = 50428
$ws.Range("N117").Value = -59606
$ws.Range("H130").Value = 89926.664
$ws.Range("J130").Value = 89926.664
$ws.Range("L130").Value = 89926.664
$ws.Range("N130").Value = -99966.664
$ws.Range("H132").Value = 20835886
$ws.Range("I132").Value = 3302.4
$ws.Range("J132").Value = 55556856
$ws.Range("K132").Value = 9907.200000000001
$ws.Range("L132").Value = 166670568
$ws.Range("M132").Value = -7377.200000000001
$ws.Range("N132").Value = -166675628

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 111111150
$ws.Range("I14").Value = 111111150
$ws.Range("K14").Value = 333333450
$ws.Range("M14").Value = -333333277
$ws.Range("H110").Value = 13617.462
$ws.Range("J110").Value = 14416.667
$ws.Range("L110").Value = 43250.001
$ws.Range("N110").Value = -51430.001
$ws.Range("H113").Value = 1482.5
$ws.Range("J113").Value = 1526.3636
$ws.Range("L113").Value = 4579.0908
$ws.Range("N113").Value = -8919.0908
$ws.Range("H122").Value = 7169.1333
$ws.Range("I122").Value = 321.63635
$ws.Range("J122").Value = 25999.75
$ws.Range("K122").Value = 2894.72715
$ws.Range("L122").Value = 233997.75
$ws.Range("M122").Value = -444.7271499999997
$ws.Range("N122").Value = -238897.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 59970.75
$ws.Range("J86").Value = 59970.75
$ws.Range("L86").Value = 59970.75
$ws.Range("N86").Value = -62342.75
$ws.Range("H89").Value = 59970.75
$ws.Range("J89").Value = 59970.75
$ws.Range("L89").Value = 179912.25
$ws.Range("N89").Value = -191768.25
$ws.Range("H99").Value = 15777.429
$ws.Range("I99").Value = 12688.4
$ws.Range("K99").Value = 12688.4
$ws.Range("M99").Value = -10442.4
$ws.Range("H113").Value = 1496.091
$ws.Range("I113").Value = 1588.8
$ws.Range("J113").Value = 1418.8334
$ws.Range("K113").Value = 1588.8
$ws.Range("L113").Value = 1418.8334
$ws.Range("M113").Value = 581.2
$ws.Range("N113").Value = -5758.8334
$ws.Range("H126").Value = 1978
$ws.Range("I126").Value = 1978
$ws.Range("K126").Value = 5934
$ws.Range("M126").Value = -3464
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("N130").Value = 0
$ws.Range("H132").Value = 3740.0833
$ws.Range("I132").Value = 3986.7778
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 11960.3334
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -9430.3334
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2750
$ws.Range("I7").Value = 2500
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -2388
$ws.Range("N7").Value = -3224
$ws.Range("H100").Value = 3386.6667
$ws.Range("I100").Value = 3080
$ws.Range("K100").Value = 3080
$ws.Range("M100").Value = -2539
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("N108").Value = 0
$ws.Range("H122").Value = 3252.8333
$ws.Range("I122").Value = 2902.4
$ws.Range("J122").Value = 5005
$ws.Range("K122").Value = 8707.200000000001
$ws.Range("L122").Value = 15015
$ws.Range("M122").Value = -6257.200000000001
$ws.Range("N122").Value = -19915
$ws.Range("H126").Value = 2750
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 3968.75
$ws.Range("I132").Value = 3850
$ws.Range("K132").Value = 11550
$ws.Range("M132").Value = -9020

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H102").Value = 33500
$ws.Range("J102").Value = 33500
$ws.Range("L102").Value = 33500
$ws.Range("N102").Value = -39990
$ws.Range("H122").Value = 2374
$ws.Range("I122").Value = 2454.375
$ws.Range("J122").Value = 2052.5
$ws.Range("K122").Value = 7363.125
$ws.Range("L122").Value = 6157.5
$ws.Range("M122").Value = -4913.125
$ws.Range("N122").Value = -11057.5
$ws.Range("H132").Value = 6784637.5
$ws.Range("I132").Value = 1620.6207
$ws.Range("J132").Value = 20835172
$ws.Range("K132").Value = 4861.8621
$ws.Range("L132").Value = 62505516
$ws.Range("M132").Value = -2331.8621
$ws.Range("N132").Value = -62510576
